$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 10, shifting existing rows 10-20 down to 11-21
$ws.Rows.Item(10).Insert()

# Update J8 (was 2 -> 3) and J9 (was 4 -> 2)
$ws.Cells.Item(8, 10).Value = 3
$ws.Cells.Item(9, 10).Value = 2

# Fill in the new row 10 with data matching the pattern of rows 8-9
$ws.Cells.Item(10, 1).Value = 121486
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = "CnC"
$ws.Cells.Item(10, 4).Value = "PICKUP_IN_STORE"
$ws.Cells.Item(10, 5).Value = "PICK"
$ws.Cells.Item(10, 6).Value = 11990
$ws.Cells.Item(10, 7).Value = 11990
$ws.Cells.Item(10, 8).Value = 11990
$ws.Cells.Item(10, 9).Value = 121486
$ws.Cells.Item(10, 10).Value = 1
$ws.Cells.Item(10, 11).Value = 457
$ws.Cells.Item(10, 12).Value = 457

# Copy the number-format style (style index 3) from row 9's A/I cells onto row 10's A/I cells
$ws.Cells.Item(9, 1).Copy()
$ws.Cells.Item(10, 1).PasteSpecial(-4122)
$ws.Cells.Item(9, 9).Copy()
$ws.Cells.Item(10, 9).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the view/selection state
$ws.Range("J10").Select()
